# Adding Test Case: Verify Lead Image #17
#
# Appends three new "lead image" test rows (Reconstruction / Reconstrucción
# article, a side-effects article, and an ibrutinib press release) to the
# pages_with_leadimage, pages_with_leadimage_Alt and pages_with_leadimage_Credit
# sheets, and the first two of those rows to pages_with_leadimage_Caption.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: pages_with_leadimage (new rows 4-6) ---------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate() | Out-Null

$ws1.Range("A4").Value = "espanol/efectos-secundarios"
$ws1.Range("F4").Value = "Reconstrucción"
$ws1.Range("D4").Value = "Los implantes se insertan debajo de la piel o del músculo del pecho"
$ws1.Range("B4").Value = "Article"
$ws1.Range("C4").Value = "Spanish"
$ws1.Range("E4").Value = "Credit: Instituto Nacional del Cáncer"

$ws1.Range("A5").Value = "about-cancer/treatment/side-effects"
$ws1.Range("D5").Value = "Tell your doctor about side effects you are experiencing, so you get the care and treatment you need to manage these problems."
$ws1.Range("F5").Value = "Reconstruction"
$ws1.Range("B5").Value = "Article"
$ws1.Range("C5").Value = "English"
$ws1.Range("E5").Value = "Credit: National Cancer Institute"

$ws1.Range("A6").Value = "news-events/press-releases/2018/leukemia-cll-ibrutinib-trial"
$ws1.Range("D6").Value = "Ibrutinib plus rituximab superior to standard treatment for some patients with chronic leukemia"
$ws1.Range("F6").Value = "patients with chronic leukemia"
$ws1.Range("B6").Value = "Press Release"
$ws1.Range("C6").Value = "English"
$ws1.Range("E6").Value = "Credit: iStock"

$ws1.Range("B20").Select() | Out-Null

# --- Sheet 3: pages_with_leadimage_Alt (new rows 7-9) -----------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate() | Out-Null

$ws3.Range("A7").Value = "espanol/efectos-secundarios"
$ws3.Range("B7").Value = "Article"
$ws3.Range("C7").Value = "Spanish"
$ws3.Range("D7").Value = "Los implantes se insertan debajo de la piel o del músculo del pecho"
$ws3.Range("E7").Value = "Credit: Instituto Nacional del Cáncer"
$ws3.Range("F7").Value = "Reconstrucción"

$ws3.Range("A8").Value = "about-cancer/treatment/side-effects"
$ws3.Range("B8").Value = "Article"
$ws3.Range("C8").Value = "English"
$ws3.Range("D8").Value = "Tell your doctor about side effects you are experiencing, so you get the care and treatment you need to manage these problems."
$ws3.Range("E8").Value = "Credit: National Cancer Institute"
$ws3.Range("F8").Value = "Reconstruction"

$ws3.Range("A9").Value = "news-events/press-releases/2018/leukemia-cll-ibrutinib-trial"
$ws3.Range("B9").Value = "Press Release"
$ws3.Range("C9").Value = "English"
$ws3.Range("D9").Value = "Ibrutinib plus rituximab superior to standard treatment for some patients with chronic leukemia"
$ws3.Range("E9").Value = "Credit: iStock"
$ws3.Range("F9").Value = "patients with chronic leukemia"

$ws3.Range("A9:XFD9").Select() | Out-Null

# --- Sheet 4: pages_with_leadimage_Credit (new rows 8-10) -------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate() | Out-Null

$ws4.Range("A8").Value = "espanol/efectos-secundarios"
$ws4.Range("B8").Value = "Article"
$ws4.Range("C8").Value = "Spanish"
$ws4.Range("D8").Value = "Los implantes se insertan debajo de la piel o del músculo del pecho"
$ws4.Range("E8").Value = "Credit: Instituto Nacional del Cáncer"
$ws4.Range("F8").Value = "Reconstrucción"

$ws4.Range("A9").Value = "about-cancer/treatment/side-effects"
$ws4.Range("B9").Value = "Article"
$ws4.Range("C9").Value = "English"
$ws4.Range("D9").Value = "Tell your doctor about side effects you are experiencing, so you get the care and treatment you need to manage these problems."
$ws4.Range("E9").Value = "Credit: National Cancer Institute"
$ws4.Range("F9").Value = "Reconstruction"

$ws4.Range("A10").Value = "news-events/press-releases/2018/leukemia-cll-ibrutinib-trial"
$ws4.Range("B10").Value = "Press Release"
$ws4.Range("C10").Value = "English"
$ws4.Range("D10").Value = "Ibrutinib plus rituximab superior to standard treatment for some patients with chronic leukemia"
$ws4.Range("E10").Value = "Credit: iStock"
$ws4.Range("F10").Value = "patients with chronic leukemia"

$ws4.Range("A10:XFD10").Select() | Out-Null

# --- Sheet 5: pages_with_leadimage_Caption (new rows 4-5) -------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate() | Out-Null

$ws5.Range("A4").Value = "espanol/efectos-secundarios"
$ws5.Range("B4").Value = "Article"
$ws5.Range("C4").Value = "Spanish"
$ws5.Range("D4").Value = "Los implantes se insertan debajo de la piel o del músculo del pecho"
$ws5.Range("E4").Value = "Credit: Instituto Nacional del Cáncer"
$ws5.Range("F4").Value = "Reconstrucción"

$ws5.Range("A5").Value = "about-cancer/treatment/side-effects"
$ws5.Range("B5").Value = "Article"
$ws5.Range("C5").Value = "English"
$ws5.Range("D5").Value = "Tell your doctor about side effects you are experiencing, so you get the care and treatment you need to manage these problems."
$ws5.Range("E5").Value = "Credit: National Cancer Institute"
$ws5.Range("F5").Value = "Reconstruction"

$ws5.Range("C11").Select() | Out-Null
